$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B3:G3")
$range.UnMerge()

$ws.Cells.Item(3, 2).Value = 42760.0393764352
$ws.Cells.Item(3, 3).Value = 42760.0393764352
$ws.Cells.Item(3, 4).Value = 42760.0393764352
$ws.Cells.Item(3, 5).Value = 42760.0393764352
$ws.Cells.Item(3, 6).Value = 42760.0393764352
$ws.Cells.Item(3, 7).Value = 42760.0393764352

$range.Merge()
